# The deck's single live theme (ppt/theme/theme2.xml, linked from
# slideMaster1.xml -> rId12) currently carries the "Integral" palette.
# The commit swaps the deck's colour scheme over to the stock
# "Office" palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink), i.e.
# exactly what Design > Variants > Colors > "Office" does in the
# PowerPoint UI. Apply it via the SlideMaster's ColorScheme, which is
# the supported COM surface for editing the active theme's 12 colour
# slots (RGBColor.RGB on each Colors(i)).

$p = $ppt.ActivePresentation
$master = $p.Designs.Item(1).SlideMaster
$scheme = $master.ColorScheme

function Set-SchemeColor($colorScheme, $index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $rgbValue = $r + ($g * 256) + ($b * 65536)
    $colorScheme.Colors($index).RGB = $rgbValue
}

# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
Set-SchemeColor $scheme 1  "000000"
Set-SchemeColor $scheme 2  "FFFFFF"
Set-SchemeColor $scheme 3  "44546A"
Set-SchemeColor $scheme 4  "E7E6E6"
Set-SchemeColor $scheme 5  "5B9BD5"
Set-SchemeColor $scheme 6  "ED7D31"
Set-SchemeColor $scheme 7  "A5A5A5"
Set-SchemeColor $scheme 8  "FFC000"
Set-SchemeColor $scheme 9  "4472C4"
Set-SchemeColor $scheme 10 "70AD47"
Set-SchemeColor $scheme 11 "0563C1"
Set-SchemeColor $scheme 12 "954F72"

Write-Output "Applied Office colour scheme to the active theme."
